$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

# --- Paragraph 24: final paragraph (HighScore struct description) rewrite ---
$xmlPara24 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Highscores</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> are saved in a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>HighScore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> struct, containing a</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> string for</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> name and </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">an integer for </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">score. When the player wants to view the existing scores, the </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>scores</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and connected names saved in the database are saved into a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>HighScore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> vector for easy access. When the player wants to save their </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>highscore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, the score is written into the database in the binary format.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(24).Range.InsertXML($xmlPara24)

# --- Paragraph 23: "Highscores" -> "Highscore struct" heading ---
$xmlPara23 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Highscore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> struct</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(23).Range.InsertXML($xmlPara23)

# --- Paragraph 22: Input paragraph rewrite, plus new blank paragraph after it ---
$xmlPara22 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The player can save their score to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>highscore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> database after each game. They can input a three-character name to save their score with, reminiscent of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>highscores</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> on old arcade machines. </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>highscore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> can only be saved if it beats an existing score in the database.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>
'@
$d.Paragraphs.Item(22).Range.InsertXML($xmlPara22)

# --- Paragraph 20: "Highscores" heading - wrap with proofErr ---
$xmlPara20 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Highscores</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$d.Paragraphs.Item(20).Range.InsertXML($xmlPara20)

# --- Paragraph 11: Shots paragraph - remove lastRenderedPageBreak ---
$xmlPara11 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The player can create shots by pressing the shoot key. These shots travel at a set speed in the direction the player was facing when they pressed the shoot key. </w:t></w:r><w:r w:rsidR="00AD6919"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Shots are initialised into a vector with a maximum size of twenty. If there are already twenty shots and the player attempts to shoot, the oldest shot will be replaced. Without a limit on how fast the player can shoot, this rewards players who time and aim their shots rather than shooting wildly as fast as they can. If a shot reaches the edge of the screen, it is removed.</w:t></w:r><w:r w:rsidR="003479CE"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> If a shot collides with an asteroid, it changes the asteroid’s velocity as if being knocked by the shot</w:t></w:r><w:r w:rsidR="00E4629F"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, which can help players escape head-on collisions with asteroids.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(11).Range.InsertXML($xmlPara11)

# --- Paragraph 9: Asteroids paragraph split into three paragraphs ---
$xmlPara9 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The asteroids are all initialised at the start of play in an array. </w:t></w:r><w:r w:rsidR="00BD32A0"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Asteroids can be small, medium, or large. This size is used to determine </w:t></w:r><w:r w:rsidR="00FB756F"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>most of</w:t></w:r><w:r w:rsidR="00BD32A0"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> the asteroid’s other characteristics, such as starting velocity, deflection from shots</w:t></w:r><w:r w:rsidR="00EB47A8"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, and health points.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">When hit by a shot, the asteroid loses a health point. If the asteroid’s health points reach zero, the player’s score is increased based on the asteroid’s size and the asteroid is respawned again after </w:t></w:r><w:r w:rsidR="001B3AAE"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>a period</w:t></w:r><w:r w:rsidR="00343F8A"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, making the game endless.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">When an asteroid is spawned, its starting position and velocity is set randomly, but there is a </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>100-unit</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> buffer based on the player’s position to stop asteroids from spawning too close to the player.</w:t></w:r><w:r w:rsidR="00B368A1"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00BD32A0"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The asteroids also endlessly loop across the screen like the Player.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(9).Range.InsertXML($xmlPara9)

# --- Paragraph 7: Player paragraph split into two paragraphs ---
$xmlPara7 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The player is the only object to have an acceleration value. This value is increased when the user holds the forward key and decrease</w:t></w:r><w:r w:rsidR="000B0181"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>d</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> when they hold backward key. </w:t></w:r><w:r w:rsidR="000B0181"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">The player’s rotation value is increased and decreased using the left and right keys. </w:t></w:r><w:r w:rsidR="00D142C6"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">By changing the values of the player’s velocity based on </w:t></w:r><w:r w:rsidR="000B0181"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">the </w:t></w:r><w:r w:rsidR="00D142C6"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>acceleration and rotation value</w:t></w:r><w:r w:rsidR="000B0181"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>s, the player can fly around the screen in any direction.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The player’s location is set to the opposite side of the screen if they pass an edge, creating an endless loop effect.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(7).Range.InsertXML($xmlPara7)

# --- Paragraph 3: Summary paragraph - wrap "Blasteroids" + split "highscore" occurrences ---
$xmlPara3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00A01B0A"><w:rPr><w:i/><w:iCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Blasteroids</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> is a successor to the classic arcade game Asteroids. It </w:t></w:r><w:r w:rsidR="00D66542"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>includes a handful of new features, such as powerups, while still maintaining the original game’s classic look and feel.</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00BD32A0"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">A ship controlled by the player flies around shooting at asteroids that randomly shoot across the screen. The aim of the game is for the player to reach a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>highscore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> before losing all their health points from colliding with asteroids.</w:t></w:r><w:r w:rsidR="00DB52C2"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> They are assisted by randomly spawning powerups that assist them</w:t></w:r><w:r w:rsidR="006D63BF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="00BD32A0"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00F62A2D"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Both the game and the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>highscore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> database are portable, allowing players to easily share their scores.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(3).Range.InsertXML($xmlPara3)

# --- Paragraph 1: Title split into two runs with proofErr wrap ---
$xmlPara1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Blasteroids</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> Design Document</w:t></w:r></w:p>

'@
$d.Paragraphs.Item(1).Range.InsertXML($xmlPara1)

Write-Output "All edits applied."
